$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest scraped figures.
# D-column values that look like plain numbers are entered with a leading
# apostrophe so Excel keeps them as text (matching the source sheet, which
# stores every Price cell as a string, e.g. "67.386.54").
$ws.Range("D2").Value = "67.386.54"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "3.220.81"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'577.80"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("D6").Value = "'182.28"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("D9").Value = "3.221.82"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("E10").Value = "  -3.20%  "
$ws.Range("E11").Value = "  -2.02%  "
$ws.Range("D12").Value = "'0.410"
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("D13").Value = "3.781.67"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").Value = "'27.66"
$ws.Range("E15").Value = "  -3.22%  "
$ws.Range("D16").Value = "67.462.36"
$ws.Range("E16").Value = "  -0.99%  "
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("D18").Value = "3.189.67"
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").Value = "'13.37"
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("D21").Value = "'394.21"
$ws.Range("E21").Value = "  +3.02%  "
$ws.Range("E22").Value = "  -2.19%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'70.80"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("E26").Value = "  -3.10%  "
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "'9.54"
$ws.Range("E28").Value = "  -2.92%  "
$ws.Range("E30").Value = "  -2.31%  "
$ws.Range("D31").Value = "'5.54"
$ws.Range("E31").Value = "  -3.27%  "
$ws.Range("D32").Value = "'22.58"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("D33").Value = "'6.97"
$ws.Range("E33").Value = "  -3.84%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  -2.50%  "
$ws.Range("D36").Value = "'161.14"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("E37").Value = "  -5.26%  "
$ws.Range("D38").Value = "'1.88"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("D42").Value = "'6.47"
$ws.Range("E42").Value = "  -4.48%  "
$ws.Range("E43").Value = "  -5.43%  "
$ws.Range("D44").Value = "'0.0682"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("D46").Value = "2.594.93"
$ws.Range("E46").Value = "  -1.96%  "
$ws.Range("E47").Value = "  -3.85%  "
$ws.Range("D48").Value = "'332.20"
$ws.Range("E48").Value = "  -4.40%  "
$ws.Range("E49").Value = "  -3.01%  "
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("E51").Value = "  -1.85%  "

# Rows 39/40: Mantle overtakes EnergySwap in ranking; refresh name, link,
# price and volume for both rows accordingly.
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").Value = "'0.803"
$ws.Range("E39").Value = "  -3.92%  "

$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'26.24"
$ws.Range("E40").Value = "  -1.61%  "
